# Add two new poll rows (id 53 = elabe 12/19, id 54 = cluster17 12/21)
# to the bottom of the data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 147: elabe poll, week 16, 12/19 ---------------------------------
$ws.Range("A147").Value = 53
$ws.Range("B147").Value = 2021
$ws.Range("C147").Value = 16
$ws.Range("D147").Value = 12
$ws.Range("E147").Value = 19
$ws.Range("F147").Value = "elabe"
$ws.Range("G147").Value = "online"
$ws.Range("H147").Value = "partially"
$ws.Range("I147").Value = 919
$ws.Range("J147").Value = 1
$ws.Range("K147").Value = 1
$ws.Range("L147").Value = 11
$ws.Range("M147").Value = 1
$ws.Range("N147").Value = 2
$ws.Range("O147").Value = 5
$ws.Range("P147").Value = 3
$ws.Range("Q147").Value = 26
$ws.Range("R147").Value = 17
$ws.Range("U147").Value = 1
$ws.Range("V147").Value = 2
$ws.Range("W147").Value = 16
$ws.Range("X147").Value = 13
$ws.Range("Y147").Value = "T_1"
$ws.Range("AA147").Value = 1
$ws.Range("AD147").Value = "T_1"

# --- Row 148: cluster17 poll, week 17, 12/21 ------------------------------
$ws.Range("A148").Value = 54
$ws.Range("B148").Value = 2021
$ws.Range("C148").Value = 17
$ws.Range("D148").Value = 12
$ws.Range("E148").Value = 21
$ws.Range("F148").Value = "cluster17"
$ws.Range("G148").Value = "online"
$ws.Range("H148").Value = "partially"
$ws.Range("I148").Value = 1419
$ws.Range("J148").Value = 1.5
$ws.Range("K148").Value = 0.5
$ws.Range("L148").Value = 12
$ws.Range("M148").Value = 1.5
$ws.Range("N148").Value = 1.5
$ws.Range("O148").Value = 5
$ws.Range("P148").Value = 2
$ws.Range("Q148").Value = 20
$ws.Range("R148").Value = 16
$ws.Range("U148").Value = 1
$ws.Range("V148").Value = 2
$ws.Range("W148").Value = 13
$ws.Range("X148").Value = 14
$ws.Range("Y148").Value = 1
$ws.Range("AA148").Value = 1
$ws.Range("AC148").Value = 7

# --- View/window bookkeeping to match the saved workbook state -----------
$ws.Range("AD148").Select()
